$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "92×77=" "89×60="
Replace-Text "49×19=" "40×26="
Replace-Text "69×45=" "86×96="
Replace-Text "72×28=" "54×39="
Replace-Text "18×97=" "78×64="
Replace-Text "36×72=" "34×17="
Replace-Text "59×57=" "50×38="
Replace-Text "35×15=" "69×77="
Replace-Text "28×41=" "47×19="
Replace-Text "86×14=" "12×49="
Replace-Text "56×16=" "74×42="
Replace-Text "82×71=" "15×69="
Replace-Text "65×79=" "93×55="
Replace-Text "97×37=" "18×28="
Replace-Text "39×15=" "50×84="
Replace-Text "33×48=" "75×78="
Replace-Text "49×30=" "36×96="
Replace-Text "48×93=" "46×21="
Replace-Text "64×90=" "27×55="
Replace-Text "61×84=" "15×57="
Replace-Text "50×14=" "50×86="
Replace-Text "77×72=" "13×98="
Replace-Text "63×68=" "39×65="
Replace-Text "22×55=" "16×86="
Replace-Text "74×65=" "79×37="
